# Adds a new worksheet "thomas" (a user configuration, like the existing
# "emre" / "tom" / "notworking" / "jan" sheets) containing a component list.

$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore the selection
# after appending the new sheet (Worksheets.Add activates the new sheet).
$activeSheetName = $wb.ActiveSheet.Name

# Insert the new sheet after the last existing sheet ("jan") so it becomes
# the final tab in the workbook, then rename it to "thomas".
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "thomas"

# Populate the two-column component table, matching the layout used by the
# other user sheets (column A = component group, column B = chosen part).
$ws.Cells.Item(1, 1).Value = "groupComponent"
$ws.Cells.Item(1, 2).Value = "nameComponent"

$ws.Cells.Item(2, 1).Value = "CPU"
$ws.Cells.Item(2, 2).Value = "i7 5960X"

$ws.Cells.Item(3, 1).Value = "RAM"
$ws.Cells.Item(3, 2).Value = "HyperX Fury"

$ws.Cells.Item(4, 1).Value = "GPU"
$ws.Cells.Item(4, 2).Value = "GTX Titan X"

$ws.Cells.Item(5, 1).Value = "PSU"
$ws.Cells.Item(5, 2).Value = "Dark Power Pro 11"

$ws.Cells.Item(6, 1).Value = "Drive"
$ws.Cells.Item(6, 2).Value = "950 EVO"

$ws.Cells.Item(7, 1).Value = "Motherboard"
$ws.Cells.Item(7, 2).Value = "MAXIMUS VII RANGER"

# Restore the previously active sheet/tab selection.
$wb.Worksheets.Item($activeSheetName).Activate()
